$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Rough title" value for the 2014 Visiting Lectureship row
$ws.Range("D13").Value = "American University in Cairo"

# Update the view: zoom scale + new selection
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("D13").Select()
